$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

foreach ($name in @("Missing Sections", "Extra Sections")) {
    $wb.Worksheets.Item($name).Delete()
}

$ws = $wb.Worksheets.Item("Summary")
$ws.Name = "Sheet1"

$headers = @(
    "Metadata Status",
    "Total ToC Entries",
    "Sections Parsed",
    "TOC Covered Pages",
    "Pages with Text",
    "Page Coverage (%)",
    "Content Coverage (%)",
    "TOC Coverage (%)",
    "JSONL Records",
    "Inheritance Detected"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Cells.Item(2, 1).Value = "Valid"
$ws.Cells.Item(2, 2).Value = 1005
$ws.Cells.Item(2, 3).Value = 1021
$ws.Cells.Item(2, 4).Value = 1031
$ws.Cells.Item(2, 5).Value = 1047
$ws.Cells.Item(2, 6).Value = 100
$ws.Cells.Item(2, 7).Value = 97.52
$ws.Cells.Item(2, 8).Value = 98.47
$ws.Cells.Item(2, 9).Value = 1021
$ws.Cells.Item(2, 10).Value = $true

# Make sure every header cell (A1:J1) has the same bold / centered / bordered
# style that the original header row (A1:D1) already used. Copy the
# formatting from the existing A1 cell onto the newly added header cells
# so they end up sharing the very same cell style.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
